# Added test for min max deductibles calcrule 13 (fm31)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# Existing tests fm29 (row 34) and fm30 (row 35) are now complete
$ws.Range("H34").Value = "complete"
$ws.Range("I34").Value = "complete"
$ws.Range("H35").Value = "complete"
$ws.Range("I35").Value = "complete"

# New row 36: fm31 - Min and Max deductibles test calcrule 13
# Copy formatting from the row above (row 35) first
$ws.Range("B35:I35").Copy()
$ws.Range("B36:I36").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B36").Value = "fm31"
$ws.Range("C36").Value = "Min and Max deductibles test calcrule 13"
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = "6, 13"
$ws.Range("F36").Value = 2
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = "in progress"
$ws.Range("I36").Value = "in progress"

# Match the new selection left by Excel after adding the row
$ws.Range("B36").Select()
